# Updated cryptos list - apply latest Price/Volume(1h) figures and the
# BitcoinCash/ShibaInu row reorder captured in the source diff.
#
# All data cells in this sheet are stored as literal text (t="inlineStr")
# even when their content looks numeric (e.g. "0.556", "33.01"). Plain
# `Range.Value = "..."` assignment lets Excel's COM layer auto-coerce
# number-looking strings into real numbers (losing the original text
# formatting / introducing float rounding noise), so price values are
# staged through a scratch cell that is forced to Text via a quote-prefix,
# then copied across with PasteSpecial-Values and cleaned up. This keeps
# the destination cell's style untouched (no NumberFormat/style churn)
# while still landing a genuine text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$stage = $ws.Range("Z1")

function Set-Text {
    param(
        [string]$Ref,
        [string]$Text
    )
    $stage.Value = "'" + $Text
    $stage.Copy()
    $ws.Range($Ref).PasteSpecial(-4163)  # xlPasteValues
}

function Set-Price {
    param(
        [int]$Row,
        [string]$Price
    )
    Set-Text "D$Row" $Price
}

function Set-Volume {
    param(
        [int]$Row,
        [string]$Pct
    )
    # Percent column values are never numeric-coercible ("  +1.23%  "),
    # so a direct assignment is safe and keeps the style untouched.
    $ws.Range("E$Row").Value = "  $Pct  "
}

# Row 2 - Bitcoin
Set-Price  2 "34.359.02"
Set-Volume 2 "+0.78%"

# Row 3 - Ethereum
Set-Price  3 "1.786.59"
Set-Volume 3 "+0.36%"

# Row 4 - TetherUSD
Set-Volume 4 "-0.05%"

# Row 5 - BNB
Set-Price  5 "226.53"
Set-Volume 5 "+0.53%"

# Row 6 - XRP
Set-Price  6 "0.556"
Set-Volume 6 "+1.96%"

# Row 7 - USDC
Set-Volume 7 "-0.07%"

# Row 8 - Solana
Set-Price  8 "33.01"
Set-Volume 8 "+3.86%"

# Row 9 - Cardano
Set-Volume 9 "+1.42%"

# Row 10 - Dogecoin
Set-Price  10 "0.0689"
Set-Volume 10 "+0.57%"

# Row 11 - TRON
Set-Price  11 "0.0947"
Set-Volume 11 "+0.04%"

# Row 12 - WrappedliquidstakedEther2.0
Set-Price  12 "2.045.87"
Set-Volume 12 "+0.42%"

# Row 13 - Chainlink
Set-Price  13 "11.25"
Set-Volume 13 "+3.21%"

# Row 14 - WrappedEther
Set-Price  14 "1.793.80"
Set-Volume 14 "+0.64%"

# Row 15 - Polygon
Set-Price  15 "0.635"
Set-Volume 15 "+2.17%"

# Row 16 - WrappedBTC
Set-Price  16 "34.398.99"
Set-Volume 16 "+0.91%"

# Row 17 - Polkadot: unchanged

# Row 18 - Litecoin
Set-Price  18 "68.46"
Set-Volume 18 "+1.27%"

# Rows 19/20 - BitcoinCash and ShibaInu swap places
Set-Text "B19" "BitcoinCash"
Set-Text "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-Price  19 "245.14"
Set-Volume 19 "+0.15%"

Set-Text "B20" "ShibaInu"
Set-Text "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-Price  20 "0.0$([char]0x2083)0795"
Set-Volume 20 "+1.30%"

# Row 21 - Avalanche
Set-Price  21 "11.31"
Set-Volume 21 "+3.52%"

# Row 22 - Dai
Set-Volume 22 "-0.06%"

# Row 23 - Uniswap
Set-Volume 23 "+1.83%"

# Row 24 - Monero
Set-Price  24 "168.41"
Set-Volume 24 "+4.33%"

# Row 25 - Toncoin
Set-Volume 25 "+1.80%"

# Row 26 - Cosmos
Set-Price  26 "7.34"
Set-Volume 26 "+3.43%"

# Row 27 - EthereumClassic
Set-Price  27 "16.54"
Set-Volume 27 "+1.94%"

# Row 28 - Stellar
Set-Volume 28 "+1.70%"

# Row 29 - BinanceUSD
Set-Volume 29 "-0.20%"

# Row 30 - InternetComputer(DFINITY)
Set-Volume 30 "+9.45%"

# Row 31 - Hedera
Set-Volume 31 "+2.02%"

# Row 32 - PancakeSwap
Set-Volume 32 "+0.83%"

# Row 33 - Filecoin
Set-Price  33 "3.80"
Set-Volume 33 "+2.63%"

# Row 34 - LidoDAOToken
Set-Volume 34 "+1.66%"

# Row 35 - Maker
Set-Price  35 "1.411.96"
Set-Volume 35 "-2.55%"

# Row 36 - RenderToken
Set-Price  36 "2.57"
Set-Volume 36 "+4.68%"

# Row 37 - ImmutableX
Set-Price  37 "0.685"
Set-Volume 37 "+5.25%"

# Row 38 - TrustWalletToken
Set-Volume 38 "+3.23%"

# Row 39 - VeChain
Set-Volume 39 "+0.38%"

# Row 40 - Aave
Set-Price  40 "84.81"
Set-Volume 40 "+5.62%"

# Row 41 - HuobiToken
Set-Volume 41 "+0.73%"

# Row 42 - ARBITRUM
Set-Volume 42 "+2.67%"

# Row 43 - MXToken
Set-Volume 43 "+2.09%"

# Row 44 - InjectiveProtocol
Set-Price  44 "13.97"
Set-Volume 44 "+2.30%"

# Row 45 - Kaspa
Set-Price  45 "0.0529"
Set-Volume 45 "+2.30%"

# Row 46 - WEMIXToken
Set-Price  46 "1.11"
Set-Volume 46 "+2.81%"

# Row 47 - FraxShare
Set-Price  47 "6.08"
Set-Volume 47 "+0.75%"

# Row 48 - RocketPoolETH
Set-Price  48 "1.946.15"
Set-Volume 48 "+0.39%"

# Row 49 - Quant
Set-Price  49 "105.41"
Set-Volume 49 "+1.26%"

# Row 50 - PaxDollar
Set-Volume 50 "-0.12%"

# Row 51 - BabyDogeCoin
Set-Volume 51 "-1.56%"

# Remove the scratch staging cell so it doesn't linger in the sheet.
$stage.Clear()
